$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("C14").Value = "параллельное звено (сингл серво)"

# Row 15
$ws.Range("C15").Value = "min"
$ws.Range("D15").Value = -90
$ws.Range("E15").Value = "middle"
$ws.Range("F15").Value = 90
$ws.Range("G15").Value = "max"

# Row 16
$ws.Range("C16").Value = 2386
$ws.Range("D16").Value = 2077.5
$ws.Range("E16").Value = 1454.5
$ws.Range("F16").Value = 822.5
$ws.Range("G16").Value = 557.5
$ws.Range("I16").Value = "mcs/ang (-90...90)"

# Row 17
$ws.Range("C17").Value = 542
$ws.Range("D17").Value = 468
$ws.Range("E17").Value = 326
$ws.Range("F17").Value = 183
$ws.Range("G17").Value = 116
$ws.Range("I17").Value = 6.972

# Row 19
$ws.Range("I19").Value = "checck"

# Row 20
$ws.Range("H20").Value = "mid+max"
$ws.Range("I20").Value = "mid+90"
$ws.Range("J20").Value = "mid-90"
$ws.Range("K20").Value = "mid-min"

# Row 21
$ws.Range("H21").Formula = "=(E16-G16)/I17"
$ws.Range("I21").Formula = "=(E16-F16)/I17"
$ws.Range("J21").Formula = "=(E16-D16)/I17"
$ws.Range("K21").Formula = "=(E16-C16)/I17"

# Row 23
$ws.Range("H23").Value = 120
$ws.Range("J23").Value = '"+120"'
$ws.Range("K23").Value = '"-120"'
$ws.Range("J23").Style = "Обычный"
$ws.Range("K23").Style = "Обычный"

# Row 24
$ws.Range("H24").Formula = "=120*I17"
$ws.Range("J24").Formula = "=E16-H24"
$ws.Range("K24").Formula = "=E16+H24"

# Apply fill style to J23, K23, J24, K24 (new style with fillId=3, theme 9 tint 0.4)
$ws.Range("J23:K24").Interior.ThemeColor = 9
$ws.Range("J23:K24").Interior.TintAndShade = -0.39997558519241921

# Selection as in diff
$ws.Range("N23").Select()
